# Apply updated cryptocurrency price/volume data to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "60.334.69"
$ws.Cells.Item(2, 5).Value = "  +0.15%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "2.336.85"
$ws.Cells.Item(3, 5).Value = "  -0.28%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  +0.01%  "

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "549.11"
$ws.Cells.Item(5, 5).Value = "  +0.60%  "

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "131.27"
$ws.Cells.Item(6, 5).Value = "  -0.78%  "

# Row 7
$ws.Cells.Item(7, 5).Value = "  +0.04%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  -1.05%  "

# Row 9
$ws.Cells.Item(9, 4).Value = "2.334.96"
$ws.Cells.Item(9, 5).Value = "  -0.22%  "

# Row 10
$ws.Cells.Item(10, 5).Value = "  +1.10%  "

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "5.62"
$ws.Cells.Item(11, 5).Value = "  +1.76%  "

# Row 12
$ws.Cells.Item(12, 5).Value = "  -0.54%  "

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "0.336"
$ws.Cells.Item(13, 5).Value = "  +0.81%  "

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "23.84"
$ws.Cells.Item(14, 5).Value = "  -0.22%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "2.753.30"
$ws.Cells.Item(15, 5).Value = "  -0.24%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "60.274.90"
$ws.Cells.Item(16, 5).Value = "  +0.12%  "

# Row 17
$ws.Cells.Item(17, 5).Value = "  +1.07%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "2.330.56"
$ws.Cells.Item(18, 5).Value = "  -0.76%  "

# Row 19
$ws.Cells.Item(19, 5).Value = "  +0.35%  "

# Row 20
$ws.Cells.Item(20, 5).Value = "  -1.37%  "

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "314.01"
$ws.Cells.Item(21, 5).Value = "  -0.05%  "

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "6.60"
$ws.Cells.Item(22, 5).Value = "  -3.26%  "

# Row 23
$ws.Cells.Item(23, 5).Value = "  +0.12%  "

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "64.11"
$ws.Cells.Item(24, 5).Value = "  +0.97%  "

# Row 25
$ws.Cells.Item(25, 5).Value = "  -1.42%  "

# Row 26
$ws.Cells.Item(26, 5).Value = "  -0.08%  "

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "7.95"
$ws.Cells.Item(27, 5).Value = "  +0.54%  "

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "1.39"
$ws.Cells.Item(28, 5).Value = "  +2.53%  "

# Row 29
$ws.Cells.Item(29, 5).Value = "  +6.94%  "

# Row 30
$ws.Cells.Item(30, 5).Value = "  -1.10%  "

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "171.35"
$ws.Cells.Item(31, 5).Value = "  -0.02%  "

# Row 32
$ws.Cells.Item(32, 4).Value = "0.0₃0735"
$ws.Cells.Item(32, 5).Value = "  +0.75%  "

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "6.08"
$ws.Cells.Item(33, 5).Value = "  +2.13%  "

# Row 34
$ws.Cells.Item(34, 2).Value = "PolygonEcosystemToken"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "0.385"
$ws.Cells.Item(34, 5).Value = "  +0.85%  "

# Row 35
$ws.Cells.Item(35, 2).Value = "ImmutableX"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "1.36"
$ws.Cells.Item(35, 5).Value = "  -3.10%  "

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "18.05"
$ws.Cells.Item(36, 5).Value = "  -0.06%  "

# Row 37
$ws.Cells.Item(37, 5).Value = "  +0.00%  "

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.999"
$ws.Cells.Item(38, 5).Value = "  -0.07%  "

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "4.12"
$ws.Cells.Item(39, 5).Value = "  -1.23%  "

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "327.58"
$ws.Cells.Item(40, 5).Value = "  +0.95%  "

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "38.14"

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "1.54"
$ws.Cells.Item(42, 5).Value = "  +0.39%  "

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "137.81"
$ws.Cells.Item(43, 5).Value = "  -2.92%  "

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "3.51"
$ws.Cells.Item(44, 5).Value = "  +1.23%  "

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.0952"
$ws.Cells.Item(45, 5).Value = "  +0.56%  "

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "19.31"
$ws.Cells.Item(46, 5).Value = "  -1.04%  "

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "0.567"
$ws.Cells.Item(47, 5).Value = "  +1.10%  "

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "0.0497"
$ws.Cells.Item(48, 5).Value = "  -0.25%  "

# Row 49
$ws.Cells.Item(49, 2).Value = "BabyDogeCoin"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Cells.Item(49, 4).Value = "0.0₆0223"
$ws.Cells.Item(49, 5).Value = "  +8.27%  "

# Row 50
$ws.Cells.Item(50, 2).Value = "VeChain"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "0.0216"
$ws.Cells.Item(50, 5).Value = "  +1.23%  "

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "10.93"
$ws.Cells.Item(51, 5).Value = "  -0.82%  "
